$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row cells: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$headerRange = $ws.Range("A1:U1")
$headerRange.Replace("_old", "_FV2404")
$headerRange.Replace("_new", "_FV2410")

# Freeze the header row (split below row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel Table with an AutoFilter
$dataRange = $ws.Range("A1:U73")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
